# Update weights (column B) and notes (column C) on Sheet1.
# Column A (resource ids R1..R22'/waste rows) is untouched - only the
# weight values and the explanatory notes text are revised, per the
# "fixed 2 initial_state files" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - R1 (population)
$ws.Range("C2").Value = "analog to population - the amount of people in a country is only a small indicator of the country's prosperity. Some small countries are very wealthy and some are very poor. Overall, does indicate some sense of wealth for a country."

# Row 3 - R2 (metallic elements)
$ws.Range("C3").Value = "analog to metallic elements; essential for metallic alloy creation and electronic creation "

# Row 4 - R3 (timber)
$ws.Range("C4").Value = "analog to timber - used in all forms of construction, but not a particularly rare resource"

# Row 5 - R4 (available land)
$ws.Range("C5").Value = "analog to available land; valued at twice the weight as water because land limits how much housing/farm/factories can be created that bring large amounts of prosperity. "

# Row 6 - R5 (renewable energy)
$ws.Range("C6").Value = "analog to renewable energy; renewable energy valued at 1 - in direct correlation to renewable energy waste's weight being -1"

# Row 7 - R6 (fossil fuel energy)
$ws.Range("B7").Value = 1.5
$ws.Range("C7").Value = "analog to fossil fuel energy - fossil fuels create more energy than green sources, but their waste is higher to indicate penalty for using nonrenewables."

# Row 8 - R7 (water)
$ws.Range("C8").Value = "analog to water; 0.5 chosen as the baseline for which all other raw resources are weighted. Essential for life and is involved in other types of resource creation, but is not rare."

# Row 9 - R8 (animals)
$ws.Range("B9").Value = 0.5
$ws.Range("C9").Value = "analog to animals; used for farms and food. Not particularly rare and has only a few use cases"

# Row 10 - R9 (plants)
$ws.Range("B10").Value = 0.5
$ws.Range("C10").Value = "analog to plants; used for farms and food, also produces fresh oxygen. Not particularly rare as well."

# Row 11 - R18 (metallic alloys)
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = "analog to metallic alloys; weighted at 2 to account for -1 alloy waste weight. Alloy + alloy waste = 2 in weight, compared to 1.5 in lost input resources"

# Row 12 - R19 (housing)
$ws.Range("B12").Value = 15
$ws.Range("C12").Value = "analog to housing; weighted at 15 to account for -2 housing waste weight. Input resources lost have combined weight of 12.25."

# Row 13 - R20 (electronics)
$ws.Range("B13").Value = 5
$ws.Range("C13").Value = "analog to electronics; weighted at 5 since 2 electronics and 1 waste is created - these total to 9 in weight compared to 8.25 of lost input resources"

# Row 17 - R5' (renewable energy waste)
$ws.Range("C17").Value = "renewable energy waste; renewable energies' waste weighted at -1 so that there is no net loss in using renewable energies"

# Row 18 - R6' (nonrenewable energy waste)
$ws.Range("B18").Value = -2
$ws.Range("C18").Value = "nonrenewable energy waste, nonrenewable energy waste is weighted higher than the weight of nonrenewable energy, to discourage fossil fuel use"

# Row 19 - R18' (metallic alloys waste)
$ws.Range("B19").Value = -1

# Row 20 - R19' (housing waste)
$ws.Range("B20").Value = -2

# Row 21 - R20' (electronics waste)
$ws.Range("B21").Value = -1

# Update the active selection to match the saved view state.
$ws.Range("H10").Select()
